# Update workbook to add data for 2021-11-09 (commit: "Add data for 2021-11-17")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab (and thus workbook.xml <sheet name=.../>)
$ws.Name = "Through 2021-11-09"

# Update the "November (through ...)" label in column A, row 12
$ws.Range("A12").Value = "November (through 11-09)"

# Update the November row (row 12) values for columns B..H (2015..2021)
$ws.Range("B12").Value = 12
$ws.Range("C12").Value = 22
$ws.Range("D12").Value = 32
$ws.Range("E12").Value = 23
$ws.Range("F12").Value = 13
$ws.Range("G12").Value = 57
$ws.Range("H12").Value = 64

# Update the Total row (row 13) values for columns B..H (2015..2021)
$ws.Range("B13").Value = 270
$ws.Range("C13").Value = 508
$ws.Range("D13").Value = 742
$ws.Range("E13").Value = 638
$ws.Range("F13").Value = 495
$ws.Range("G13").Value = 1114
$ws.Range("H13").Value = 1508
